$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data rows -----------------------------------------------------------
# Row 2: Combinations / Recursion (note leading space before "Recursion")
$ws.Range("A2").Value = "Combinations"
$ws.Range("B2").Value = " Recursion"

# Row 3: Subsets / Recursion
$ws.Range("A3").Value = "Subsets"
$ws.Range("B3").Value = "Recursion"

# Row 4: Subsets II / Recursion
$ws.Range("A4").Value = "Subsets II"
$ws.Range("B4").Value = "Recursion"

# Row 5: Generate Parentheses / Recursion
$ws.Range("A5").Value = "Generate Parentheses"
$ws.Range("B5").Value = "Recursion"

# --- Hyperlinks on column A -----------------------------------------------
$ws.Hyperlinks.Add($ws.Range("A2"), "https://leetcode.com/problems/combinations/")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://leetcode.com/problems/subsets/")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://leetcode.com/problems/subsets-ii/")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://leetcode.com/problems/generate-parentheses/")

# --- Formatting ------------------------------------------------------------
# Column A links (A2 & A5 plain left/top default, A3 & A4 left+top aligned)
$linkRange = $ws.Range("A2:A5")
$linkRange.Style = "Normal"
$linkRange.Font.Name = "Arial"
$linkRange.Font.Size = 10
$linkRange.Font.Underline = $true
$linkRange.Font.Color = 16711680

$ws.Range("A3:A4").HorizontalAlignment = -4131
$ws.Range("A3:A4").VerticalAlignment = -4160

# Column B category cells
$catRange = $ws.Range("B2:B5")
$catRange.Font.Name = "Aptos Narrow"
$catRange.Font.Size = 10
$catRange.HorizontalAlignment = -4131
$catRange.Interior.Pattern = 1
$catRange.Interior.Color = 16777215
$catRange.Interior.PatternColor = 16777215

$ws.Range("B5").Select
